{"js": "// Replace the multiplication-problem text in the table cells per the\n// authoring diff. Each old expression is unique in the document, so we\n// look each one up by an exact, case-sensitive search and replace just\n// that run's text, preserving all existing run/paragraph formatting.\nconst replacements = [\n  [\"659\u00d79=\", \"670\u00d76=\"],\n  [\"533\u00d73=\", \"394\u00d77=\"],\n  [\"942\u00d72=\", \"268\u00d73=\"],\n  [\"578\u00d77=\", \"704\u00d76=\"],\n  [\"939\u00d78=\", \"468\u00d75=\"],\n  [\"373\u00d75=\", \"989\u00d77=\"],\n  [\"815\u00d77=\", \"579\u00d79=\"],\n  [\"803\u00d78=\", \"946\u00d75=\"],\n  [\"441\u00d72=\", \"882\u00d75=\"],\n  [\"149\u00d73=\", \"652\u00d73=\"],\n  [\"972\u00d75=\", \"147\u00d76=\"],\n  [\"436\u00d79=\", \"939\u00d78=\"],\n  [\"583\u00d79=\", \"680\u00d72=\"],\n  [\"569\u00d78=\", \"284\u00d76=\"],\n  [\"277\u00d76=\", \"176\u00d73=\"],\n  [\"267\u00d73=\", \"441\u00d76=\"],\n  [\"525\u00d77=\", \"675\u00d76=\"],\n  [\"459\u00d79=\", \"346\u00d76=\"],\n  [\"112\u00d79=\", \"990\u00d78=\"],\n  [\"556\u00d72=\", \"334\u00d76=\"],\n  [\"253\u00d76=\", \"931\u00d72=\"],\n  [\"232\u00d78=\", \"772\u00d77=\"],\n  [\"440\u00d74=\", \"422\u00d78=\"],\n  [\"492\u00d75=\", \"828\u00d76=\"],\n  [\"621\u00d79=\", \"652\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in the table cells per the\n# authoring diff. Each old expression is unique in the document, so a\n# simple Find/Replace (wdReplaceAll, but only ever one hit) for each\n# exact old string is sufficient and keeps existing run/paragraph\n# formatting untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"659\u00d79=\", \"670\u00d76=\"),\n    @(\"533\u00d73=\", \"394\u00d77=\"),\n    @(\"942\u00d72=\", \"268\u00d73=\"),\n    @(\"578\u00d77=\", \"704\u00d76=\"),\n    @(\"939\u00d78=\", \"468\u00d75=\"),\n    @(\"373\u00d75=\", \"989\u00d77=\"),\n    @(\"815\u00d77=\", \"579\u00d79=\"),\n    @(\"803\u00d78=\", \"946\u00d75=\"),\n    @(\"441\u00d72=\", \"882\u00d75=\"),\n    @(\"149\u00d73=\", \"652\u00d73=\"),\n    @(\"972\u00d75=\", \"147\u00d76=\"),\n    @(\"436\u00d79=\", \"939\u00d78=\"),\n    @(\"583\u00d79=\", \"680\u00d72=\"),\n    @(\"569\u00d78=\", \"284\u00d76=\"),\n    @(\"277\u00d76=\", \"176\u00d73=\"),\n    @(\"267\u00d73=\", \"441\u00d76=\"),\n    @(\"525\u00d77=\", \"675\u00d76=\"),\n    @(\"459\u00d79=\", \"346\u00d76=\"),\n    @(\"112\u00d79=\", \"990\u00d78=\"),\n    @(\"556\u00d72=\", \"334\u00d76=\"),\n    @(\"253\u00d76=\", \"931\u00d72=\"),\n    @(\"232\u00d78=\", \"772\u00d77=\"),\n    @(\"440\u00d74=\", \"422\u00d78=\"),\n    @(\"492\u00d75=\", \"828\u00d76=\"),\n    @(\"621\u00d79=\", \"652\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
